$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 17-22: copy column-A style (bold/border/center) from A16 template ---
$ws.Cells.Item(16, 1).Copy($ws.Range("A17:A22"))

# Row 7: flowbotallopen / raw
$ws.Cells.Item(7, 2).Value = 'flowbotallopen'
$ws.Cells.Item(7, 3).Value = 'raw'
$ws.Cells.Item(7, 4).Value = 0
$ws.Cells.Item(7, 5).Value = 0.5775957040770099
$ws.Cells.Item(7, 6).Value = 0.5872966996688901
$ws.Cells.Item(7, 7).Value = 0.9435511942718612
$ws.Cells.Item(7, 8).Value = 0.726889918262501
$ws.Cells.Item(7, 9).Value = 1
$ws.Cells.Item(7, 10).Value = 0.4581593895024392
$ws.Cells.Item(7, 11).Value = 0.8674649961167815
$ws.Cells.Item(7, 12).Value = 0.8420231239051486
$ws.Cells.Item(7, 13).Value = 0.02006576500047157
$ws.Cells.Item(7, 14).Value = 0.6414610853430778
$ws.Cells.Item(7, 15).Value = 0.8644771565725876
$ws.Cells.Item(7, 16).Value = 0.7747019942539959
$ws.Cells.Item(7, 17).Value = 0.7308324085738898
$ws.Cells.Item(7, 18).Value = 0.7975350667878882
$ws.Cells.Item(7, 19).Value = 0.9218991476965369
$ws.Cells.Item(7, 20).Value = 0.0002909031672504
$ws.Cells.Item(7, 21).Value = 0.5747286531041927
$ws.Cells.Item(7, 22).Value = 0.921810695179186
$ws.Cells.Item(7, 23).Value = 0.8962025627950905
$ws.Cells.Item(7, 24).Value = 0.90144058608341

# Row 8: flowbot / sgp
$ws.Cells.Item(8, 2).Value = 'flowbot'
$ws.Cells.Item(8, 3).Value = 'sgp'
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = 0.7138278559596645
$ws.Cells.Item(8, 6).Value = 0.9289711072750856
$ws.Cells.Item(8, 7).Value = 0.9514460108732092
$ws.Cells.Item(8, 8).Value = 0.9518634328757024
$ws.Cells.Item(8, 9).Value = 1
$ws.Cells.Item(8, 10).Value = 0.9370209920996716
$ws.Cells.Item(8, 11).Value = 0.9242445997784808
$ws.Cells.Item(8, 12).Value = 0.8589363510185869
$ws.Cells.Item(8, 13).Value = 0.3082338336756623
$ws.Cells.Item(8, 14).Value = 0.933658158614168
$ws.Cells.Item(8, 15).Value = 0.9514086154617384
$ws.Cells.Item(8, 16).Value = 0.8948717248922784
$ws.Cells.Item(8, 17).Value = 0.8518262214254619
$ws.Cells.Item(8, 18).Value = 0.9263976798275305
$ws.Cells.Item(8, 19).Value = 0.9149756229122948
$ws.Cells.Item(8, 20).Value = 0.9257741993657416
$ws.Cells.Item(8, 21).Value = 0.3157468635106495
$ws.Cells.Item(8, 22).Value = 0.9262509398520498
$ws.Cells.Item(8, 23).Value = 0.9637726579785376
$ws.Cells.Item(8, 24).Value = 0.9240362455077578

# Row 9: dit / sgp
$ws.Cells.Item(9, 2).Value = 'dit'
$ws.Cells.Item(9, 3).Value = 'sgp'
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(9, 5).Value = 0.4888263669975369
$ws.Cells.Item(9, 6).Value = 0.4285037293703598
$ws.Cells.Item(9, 7).Value = 0.5999187990409752
$ws.Cells.Item(9, 8).Value = 0.9409981636998286
$ws.Cells.Item(9, 9).Value = 0.9993729480627284
$ws.Cells.Item(9, 10).Value = 0.9345694197791012
$ws.Cells.Item(9, 11).Value = 0.9404691822006718
$ws.Cells.Item(9, 12).Value = 0.8490230578053026
$ws.Cells.Item(9, 13).Value = 0.3057283838462095
$ws.Cells.Item(9, 14).Value = 0.931907603190636
$ws.Cells.Item(9, 15).Value = 0.9045179583939842
$ws.Cells.Item(9, 16).Value = 0.8692535253098979
$ws.Cells.Item(9, 17).Value = 0.7526504346555315
$ws.Cells.Item(9, 18).Value = 0.8555069896328563
$ws.Cells.Item(9, 19).Value = 0.4702669908623247
$ws.Cells.Item(9, 20).Value = 0.2080094774270662
$ws.Cells.Item(9, 21).Value = 0.6296906794596634
$ws.Cells.Item(9, 22).Value = 0.9405993364337688
$ws.Cells.Item(9, 23).Value = 0.7484414967927199
$ws.Cells.Item(9, 24).Value = 0.9185602400683186

# Row 10: pndit / sgp
$ws.Cells.Item(10, 2).Value = 'pndit'
$ws.Cells.Item(10, 3).Value = 'sgp'
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = 0.6629126197714845
$ws.Cells.Item(10, 6).Value = 0.7466600982420376
$ws.Cells.Item(10, 7).Value = 0.8899383246470043
$ws.Cells.Item(10, 8).Value = 0.7415691147898343
$ws.Cells.Item(10, 9).Value = 1
$ws.Cells.Item(10, 10).Value = 0.9487199560367378
$ws.Cells.Item(10, 11).Value = 0.9108832142439263
$ws.Cells.Item(10, 12).Value = 0.8889656155154814
$ws.Cells.Item(10, 13).Value = 0.8401693753933762
$ws.Cells.Item(10, 14).Value = 0.8325445608477199
$ws.Cells.Item(10, 15).Value = 0.9144821807549848
$ws.Cells.Item(10, 16).Value = 0.8824768579847121
$ws.Cells.Item(10, 17).Value = 0.9044884399165648
$ws.Cells.Item(10, 18).Value = 0.844398905571329
$ws.Cells.Item(10, 19).Value = 0.5970163049040024
$ws.Cells.Item(10, 20).Value = 0.6507120977643042
$ws.Cells.Item(10, 21).Value = 0.3344129876238517
$ws.Cells.Item(10, 22).Value = 0.938724984809892
$ws.Cells.Item(10, 23).Value = 0.9602763930918902
$ws.Cells.Item(10, 24).Value = 0.9332042963984388

# Row 11: hisdit / sgp
$ws.Cells.Item(11, 2).Value = 'hisdit'
$ws.Cells.Item(11, 3).Value = 'sgp'
$ws.Cells.Item(11, 4).Value = 0
$ws.Cells.Item(11, 5).Value = 0.7308486472314362
$ws.Cells.Item(11, 6).Value = 0.5122527009758682
$ws.Cells.Item(11, 7).Value = 0.8247648859714521
$ws.Cells.Item(11, 8).Value = 0.9307432024874372
$ws.Cells.Item(11, 9).Value = 1
$ws.Cells.Item(11, 10).Value = 0.9424919659479258
$ws.Cells.Item(11, 11).Value = 0.9539296841276552
$ws.Cells.Item(11, 12).Value = 0.6
$ws.Cells.Item(11, 13).Value = 0.9642704773391272
$ws.Cells.Item(11, 14).Value = 0.9313987455565628
$ws.Cells.Item(11, 15).Value = 0.9096867513268464
$ws.Cells.Item(11, 16).Value = 0.921167023838839
$ws.Cells.Item(11, 17).Value = 0.8219552328814616
$ws.Cells.Item(11, 18).Value = 0.7751844694860945
$ws.Cells.Item(11, 19).Value = 0.931173282872618
$ws.Cells.Item(11, 20).Value = 0.7614469929052529
$ws.Cells.Item(11, 21).Value = 0.6316867420732576
$ws.Cells.Item(11, 22).Value = 0.9430531400964764
$ws.Cells.Item(11, 23).Value = 0.8377400911706901
$ws.Cells.Item(11, 24).Value = 0.9308713802145232

# Row 12: pnhisdit / sgp
$ws.Cells.Item(12, 2).Value = 'pnhisdit'
$ws.Cells.Item(12, 3).Value = 'sgp'
$ws.Cells.Item(12, 4).Value = 0
$ws.Cells.Item(12, 5).Value = 0.7064115310764051
$ws.Cells.Item(12, 6).Value = 0.5139790013335145
$ws.Cells.Item(12, 7).Value = 0.9253996850275632
$ws.Cells.Item(12, 8).Value = 0.774235605073848
$ws.Cells.Item(12, 9).Value = 0.9995528758264598
$ws.Cells.Item(12, 10).Value = 0.9605903252381828
$ws.Cells.Item(12, 11).Value = 0.9339458920182081
$ws.Cells.Item(12, 12).Value = 0.8079896561887102
$ws.Cells.Item(12, 13).Value = 0.3101927743762022
$ws.Cells.Item(12, 14).Value = 0.937363334411392
$ws.Cells.Item(12, 15).Value = 0.9631899269363112
$ws.Cells.Item(12, 16).Value = 0.9329073880727108
$ws.Cells.Item(12, 17).Value = 0.939052657691387
$ws.Cells.Item(12, 18).Value = 0.946467530764085
$ws.Cells.Item(12, 19).Value = 0.934873387371325
$ws.Cells.Item(12, 20).Value = 0.8968881266656368
$ws.Cells.Item(12, 21).Value = 0.817796311628736
$ws.Cells.Item(12, 22).Value = 0.9360559140643132
$ws.Cells.Item(12, 23).Value = 0.8537877913754456
$ws.Cells.Item(12, 24).Value = 0.9328740782906444

# Row 13: hisditonly / sgp
$ws.Cells.Item(13, 2).Value = 'hisditonly'
$ws.Cells.Item(13, 3).Value = 'sgp'
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0.961556856628588
$ws.Cells.Item(13, 6).Value = 0.0889405835347683
$ws.Cells.Item(13, 7).Value = 0.2613828534507007
$ws.Cells.Item(13, 8).Value = 0.9410201566134628
$ws.Cells.Item(13, 9).Value = 1
$ws.Cells.Item(13, 10).Value = 0.9396101939917234
$ws.Cells.Item(13, 11).Value = 0.9589084133678036
$ws.Cells.Item(13, 12).Value = 0.8
$ws.Cells.Item(13, 13).Value = 0.6350873887600679
$ws.Cells.Item(13, 14).Value = 0.7902257641030861
$ws.Cells.Item(13, 15).Value = 0.932266467500007
$ws.Cells.Item(13, 16).Value = 0.8928541108167063
$ws.Cells.Item(13, 17).Value = 0.9037870681214863
$ws.Cells.Item(13, 18).Value = 0.9082900596185676
$ws.Cells.Item(13, 19).Value = 0.9229921306942644
$ws.Cells.Item(13, 20).Value = 0.0226352728866066
$ws.Cells.Item(13, 21).Value = 0.3537579683589634
$ws.Cells.Item(13, 22).Value = 0.93770971462465
$ws.Cells.Item(13, 23).Value = 0.6502681884210223
$ws.Cells.Item(13, 24).Value = 0.5771985225235223

# Row 14: hisonly / sgp
$ws.Cells.Item(14, 2).Value = 'hisonly'
$ws.Cells.Item(14, 3).Value = 'sgp'
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 0.5544635472149846
$ws.Cells.Item(14, 6).Value = 0.2169871603250254
$ws.Cells.Item(14, 7).Value = 0.25
$ws.Cells.Item(14, 8).Value = 0.9589802279644588
$ws.Cells.Item(14, 9).Value = 1
$ws.Cells.Item(14, 10).Value = 0.4747361688593852
$ws.Cells.Item(14, 11).Value = 0.9484331168438956
$ws.Cells.Item(14, 12).Value = 0.8
$ws.Cells.Item(14, 13).Value = 0.02981696275251493
$ws.Cells.Item(14, 14).Value = 0.7582339347824775
$ws.Cells.Item(14, 15).Value = 0.9211910785569832
$ws.Cells.Item(14, 16).Value = 0.8453338233183735
$ws.Cells.Item(14, 17).Value = 0.7943689685785805
$ws.Cells.Item(14, 18).Value = 0.7508333214522451
$ws.Cells.Item(14, 19).Value = 0.8102698814761254
$ws.Cells.Item(14, 20).Value = 0.0118553015307967
$ws.Cells.Item(14, 21).Value = 0.3062005729272828
$ws.Cells.Item(14, 22).Value = 0.8153569488872587
$ws.Cells.Item(14, 23).Value = 0.6746233693233445
$ws.Cells.Item(14, 24).Value = 0.800184658812862

# Row 15: pndit&pn++ / sgp
$ws.Cells.Item(15, 2).Value = 'pndit&pn++'
$ws.Cells.Item(15, 3).Value = 'sgp'
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0.4691307207856198
$ws.Cells.Item(15, 6).Value = 0.9400283511907082
$ws.Cells.Item(15, 7).Value = 0.7322693818116437
$ws.Cells.Item(15, 8).Value = 0.8508917738562736
$ws.Cells.Item(15, 9).Value = 1
$ws.Cells.Item(15, 10).Value = 0.925228467493984
$ws.Cells.Item(15, 11).Value = 0.956889624458642
$ws.Cells.Item(15, 12).Value = 0.8900924665670253
$ws.Cells.Item(15, 13).Value = 0.3063550206833829
$ws.Cells.Item(15, 14).Value = 0.9348710260721662
$ws.Cells.Item(15, 15).Value = 0.9363529922657632
$ws.Cells.Item(15, 16).Value = 0.914531368814017
$ws.Cells.Item(15, 17).Value = 0.9276253228939298
$ws.Cells.Item(15, 18).Value = 0.932984219744539
$ws.Cells.Item(15, 19).Value = 0.9449465116806404
$ws.Cells.Item(15, 20).Value = 0.9206073802550248
$ws.Cells.Item(15, 21).Value = 0.3584884745848697
$ws.Cells.Item(15, 22).Value = 0.9354923638755008
$ws.Cells.Item(15, 23).Value = 0.8717896628163613
$ws.Cells.Item(15, 24).Value = 0.9346662340419836

# Row 16: dit&pn++ / sgp
$ws.Cells.Item(16, 2).Value = 'dit&pn++'
$ws.Cells.Item(16, 3).Value = 'sgp'
$ws.Cells.Item(16, 4).Value = 0
$ws.Cells.Item(16, 5).Value = 0.7301867218910854
$ws.Cells.Item(16, 6).Value = 0.9317058211891958
$ws.Cells.Item(16, 7).Value = 0.8627579999030595
$ws.Cells.Item(16, 8).Value = 0.9479817257937072
$ws.Cells.Item(16, 9).Value = 1
$ws.Cells.Item(16, 10).Value = 0.918611770065529
$ws.Cells.Item(16, 11).Value = 0.9771820886466894
$ws.Cells.Item(16, 12).Value = 0.8377455714585403
$ws.Cells.Item(16, 13).Value = 0.3129463861646694
$ws.Cells.Item(16, 14).Value = 0.934907087247616
$ws.Cells.Item(16, 15).Value = 0.9179893990231394
$ws.Cells.Item(16, 16).Value = 0.9093645087721356
$ws.Cells.Item(16, 17).Value = 0.8243519357848043
$ws.Cells.Item(16, 18).Value = 0.9428792186473928
$ws.Cells.Item(16, 19).Value = 0.9293936465152882
$ws.Cells.Item(16, 20).Value = 0.9228225167519656
$ws.Cells.Item(16, 21).Value = 0.3584748582902033
$ws.Cells.Item(16, 22).Value = 0.929701927956842
$ws.Cells.Item(16, 23).Value = 0.8753540034520791
$ws.Cells.Item(16, 24).Value = 0.910639404277714

# Row 17: largedit&pn++ / sgp
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = 'largedit&pn++'
$ws.Cells.Item(17, 3).Value = 'sgp'
$ws.Cells.Item(17, 4).Value = 0
$ws.Cells.Item(17, 5).Value = 0.7273269810970115
$ws.Cells.Item(17, 6).Value = 0.8995981881489632
$ws.Cells.Item(17, 7).Value = 0.7109559285660698
$ws.Cells.Item(17, 8).Value = 0.9395127144785804
$ws.Cells.Item(17, 9).Value = 1
$ws.Cells.Item(17, 10).Value = 0.9515306872093762
$ws.Cells.Item(17, 11).Value = 0.927219227743634
$ws.Cells.Item(17, 12).Value = 0.8449673620617937
$ws.Cells.Item(17, 13).Value = 0.3053521311551217
$ws.Cells.Item(17, 14).Value = 0.9268652402953752
$ws.Cells.Item(17, 15).Value = 0.8053098102033086
$ws.Cells.Item(17, 16).Value = 0.9034545771158181
$ws.Cells.Item(17, 17).Value = 0.7229976504790084
$ws.Cells.Item(17, 18).Value = 0.7579144059760083
$ws.Cells.Item(17, 19).Value = 0.9415611758402678
$ws.Cells.Item(17, 20).Value = 0.8875957653621757
$ws.Cells.Item(17, 21).Value = 0.956898775000406
$ws.Cells.Item(17, 22).Value = 0.9301295195484586
$ws.Cells.Item(17, 23).Value = 0.8866346314810067
$ws.Cells.Item(17, 24).Value = 0.9329968447669564

# Row 18: flowbotallopen / sgp
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = 'flowbotallopen'
$ws.Cells.Item(18, 3).Value = 'sgp'
$ws.Cells.Item(18, 4).Value = 0
$ws.Cells.Item(18, 5).Value = 0.6810513221671851
$ws.Cells.Item(18, 6).Value = 0.6417785242810977
$ws.Cells.Item(18, 7).Value = 0.9273562855000586
$ws.Cells.Item(18, 8).Value = 0.9576042138624502
$ws.Cells.Item(18, 9).Value = 1
$ws.Cells.Item(18, 10).Value = 0.7807644094262043
$ws.Cells.Item(18, 11).Value = 0.8711505701242831
$ws.Cells.Item(18, 12).Value = 0.8434000822207519
$ws.Cells.Item(18, 13).Value = 0.3333333333333333
$ws.Cells.Item(18, 14).Value = 0.79928164814905
$ws.Cells.Item(18, 15).Value = 0.9544092831270539
$ws.Cells.Item(18, 16).Value = 0.8959796651393381
$ws.Cells.Item(18, 17).Value = 0.8360227341682046
$ws.Cells.Item(18, 18).Value = 0.9572299156137444
$ws.Cells.Item(18, 19).Value = 0.9303038269138458
$ws.Cells.Item(18, 20).Value = 0.0149991764431838
$ws.Cells.Item(18, 21).Value = 0.3853901168825424
$ws.Cells.Item(18, 22).Value = 0.923231732317806
$ws.Cells.Item(18, 23).Value = 0.8475350762788363
$ws.Cells.Item(18, 24).Value = 0.8574471090636353

# Row 19: pndit&pn++allopen>005 / sgp
$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = 'pndit&pn++allopen>005'
$ws.Cells.Item(19, 3).Value = 'sgp'
$ws.Cells.Item(19, 4).Value = 0
$ws.Cells.Item(19, 5).Value = 0.6296061724169432
$ws.Cells.Item(19, 6).Value = 0.922255040912476
$ws.Cells.Item(19, 7).Value = 0.9390989833054656
$ws.Cells.Item(19, 8).Value = 0.9676421289570594
$ws.Cells.Item(19, 9).Value = 1
$ws.Cells.Item(19, 10).Value = 0.950465941929064
$ws.Cells.Item(19, 11).Value = 0.9222544849126231
$ws.Cells.Item(19, 12).Value = 0.837403060953738
$ws.Cells.Item(19, 13).Value = 0.4368373423512981
$ws.Cells.Item(19, 14).Value = 0.9338466924047708
$ws.Cells.Item(19, 15).Value = 0.9054860799661724
$ws.Cells.Item(19, 16).Value = 0.943067027671105
$ws.Cells.Item(19, 17).Value = 0.9389742083057596
$ws.Cells.Item(19, 18).Value = 0.7586843838535584
$ws.Cells.Item(19, 19).Value = 0.6313566429005669
$ws.Cells.Item(19, 20).Value = 0.921188757063587
$ws.Cells.Item(19, 21).Value = 0.3814672994790609
$ws.Cells.Item(19, 22).Value = 0.9396433026330286
$ws.Cells.Item(19, 23).Value = 0.8677833542640851
$ws.Cells.Item(19, 24).Value = 0.9404987523568654

# Row 20: dit&pn++allopen>01 / sgp
$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 2).Value = 'dit&pn++allopen>01'
$ws.Cells.Item(20, 3).Value = 'sgp'
$ws.Cells.Item(20, 4).Value = 0
$ws.Cells.Item(20, 5).Value = 0.7189512005241321
$ws.Cells.Item(20, 6).Value = 0.9485654653040948
$ws.Cells.Item(20, 7).Value = 0.7041576718646176
$ws.Cells.Item(20, 8).Value = 0.8503720907971369
$ws.Cells.Item(20, 9).Value = 0.8
$ws.Cells.Item(20, 10).Value = 0.9449888747185105
$ws.Cells.Item(20, 11).Value = 0.9125739998881752
$ws.Cells.Item(20, 12).Value = 0.6423432299286097
$ws.Cells.Item(20, 13).Value = 0.3099760283648656
$ws.Cells.Item(20, 14).Value = 0.7899385039863834
$ws.Cells.Item(20, 15).Value = 0.8609342291449414
$ws.Cells.Item(20, 16).Value = 0.8887371525787127
$ws.Cells.Item(20, 17).Value = 0.8352117209873667
$ws.Cells.Item(20, 18).Value = 0.756063711584696
$ws.Cells.Item(20, 19).Value = 0.6333896704647582
$ws.Cells.Item(20, 20).Value = 0.9317831109813114
$ws.Cells.Item(20, 21).Value = 0.6464667381068901
$ws.Cells.Item(20, 22).Value = 0.9373351108470954
$ws.Cells.Item(20, 23).Value = 0.8268744124933202
$ws.Cells.Item(20, 24).Value = 0.6941901958087614

# Row 21: hispnditckpt299 / sgp
$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(21, 2).Value = 'hispnditckpt299'
$ws.Cells.Item(21, 3).Value = 'sgp'
$ws.Cells.Item(21, 4).Value = 0
$ws.Cells.Item(21, 5).Value = 0.6383061445139828
$ws.Cells.Item(21, 6).Value = 0.4859589221426986
$ws.Cells.Item(21, 7).Value = 0.8309260976759134
$ws.Cells.Item(21, 8).Value = 0.9559541569516719
$ws.Cells.Item(21, 9).Value = 1
$ws.Cells.Item(21, 10).Value = 0.9748837861596176
$ws.Cells.Item(21, 11).Value = 0.9389391153707609
$ws.Cells.Item(21, 12).Value = 0.8933586607942084
$ws.Cells.Item(21, 13).Value = 0.1934708289242569
$ws.Cells.Item(21, 14).Value = 0.937902792284849
$ws.Cells.Item(21, 15).Value = 0.9485533316394013
$ws.Cells.Item(21, 16).Value = 0.9208433663887449
$ws.Cells.Item(21, 17).Value = 0.9187426819475297
$ws.Cells.Item(21, 18).Value = 0.7213693535465666
$ws.Cells.Item(21, 19).Value = 0.8691632039796019
$ws.Cells.Item(21, 20).Value = 0.9313143224087592
$ws.Cells.Item(21, 21).Value = 0.314294727323924
$ws.Cells.Item(21, 22).Value = 0.9242476046829284
$ws.Cells.Item(21, 23).Value = 0.9025091162043726
$ws.Cells.Item(21, 24).Value = 0.9121223877874056

# Row 22: None / None
$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(22, 5).Value = 4
$ws.Cells.Item(22, 6).Value = 2
$ws.Cells.Item(22, 7).Value = 4
$ws.Cells.Item(22, 8).Value = 11
$ws.Cells.Item(22, 9).Value = 5
$ws.Cells.Item(22, 10).Value = 2
$ws.Cells.Item(22, 11).Value = 11
$ws.Cells.Item(22, 12).Value = 5
$ws.Cells.Item(22, 13).Value = 3
$ws.Cells.Item(22, 14).Value = 15
$ws.Cells.Item(22, 15).Value = 49
$ws.Cells.Item(22, 16).Value = 153
$ws.Cells.Item(22, 17).Value = 9
$ws.Cells.Item(22, 18).Value = 5
$ws.Cells.Item(22, 19).Value = 3
$ws.Cells.Item(22, 20).Value = 1
$ws.Cells.Item(22, 21).Value = 3
$ws.Cells.Item(22, 22).Value = 8
$ws.Cells.Item(22, 23).Value = 7
$ws.Cells.Item(22, 24).Value = 4
